$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.62780766666667
$ws.Range("H2").Value = 280.883423
$ws.Range("I2").Value = 0.3228593149748609
$ws.Range("J2").Value = 0.3228593149748609
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 2.180217129326
$ws.Range("R2").Value = 19.621954163934
$ws.Range("S2").Value = 0.003006049606144582
$ws.Range("T2").Value = 0.003006049606144583
$ws.Range("G3").Value = 93.62780766666667
$ws.Range("H3").Value = 280.883423
$ws.Range("I3").Value = 0.3228593149748609
$ws.Range("J3").Value = 0.3228593149748609
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 12.44494577651489
$ws.Range("R3").Value = 112.004511988634
$ws.Range("S3").Value = 0.01715889846326836
$ws.Range("T3").Value = 0.01715889846326836
$ws.Range("G4").Value = 93.62780766666667
$ws.Range("H4").Value = 280.883423
$ws.Range("I4").Value = 0.3228593149748609
$ws.Range("J4").Value = 0.3228593149748609
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 219.5371102089543
$ws.Range("R4").Value = 1975.833991880588
$ws.Range("S4").Value = 0.3026943669054479
$ws.Range("T4").Value = 0.3026943669054479
$ws.Range("G5").Value = 66.39541
$ws.Range("I5").Value = 0.228953097635189
$ws.Range("J5").Value = 0.228953097635189
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 1.54608351726
$ws.Range("R5").Value = 13.91475165534
$ws.Range("S5").Value = 0.002131716004617774
$ws.Range("T5").Value = 0.002131716004617774
$ws.Range("G6").Value = 66.39541
$ws.Range("I6").Value = 0.228953097635189
$ws.Range("J6").Value = 0.228953097635189
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("Q6").Value = 8.825233633593331
$ws.Range("S6").Value = 0.01216809543029251
$ws.Range("T6").Value = 0.01216809543029251
$ws.Range("G7").Value = 66.39541
$ws.Range("I7").Value = 0.228953097635189
$ws.Range("J7").Value = 0.228953097635189
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 155.6829835686534
$ws.Range("R7").Value = 1401.14685211788
$ws.Range("S7").Value = 0.2146532862002787
$ws.Range("T7").Value = 0.2146532862002787
$ws.Range("G8").Value = 129.9724656666667
$ws.Range("H8").Value = 389.917397
$ws.Range("I8").Value = 0.4481875873899502
$ws.Range("J8").Value = 0.4481875873899502
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 3.026538835514
$ws.Range("R8").Value = 27.238849519626
$ws.Range("S8").Value = 0.004172944865033101
$ws.Range("T8").Value = 0.004172944865033102
$ws.Range("G9").Value = 129.9724656666667
$ws.Range("H9").Value = 389.917397
$ws.Range("I9").Value = 0.4481875873899502
$ws.Range("J9").Value = 0.4481875873899502
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 17.27585348810289
$ws.Range("R9").Value = 155.482681392926
$ws.Range("S9").Value = 0.02381967918478727
$ws.Range("T9").Value = 0.02381967918478727
$ws.Range("G10").Value = 129.9724656666667
$ws.Range("H10").Value = 389.917397
$ws.Range("I10").Value = 0.4481875873899502
$ws.Range("J10").Value = 0.4481875873899502
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 304.7575312323702
$ws.Range("R10").Value = 2742.817781091332
$ws.Range("S10").Value = 0.4201949633401299
$ws.Range("T10").Value = 0.4201949633401299
